# DEBIT_BOOK.xlsx edit:
#  - Insert a new row at 18 in the right-hand mini table (date 45293 / 10)
#  - Bump the "Duy chuyen cho co Diem" transfer on the last existing entry
#    from 700k to 1.7tr (i.e. +1,000k = +1tr)
#  - Add a new final entry recording a 10tr loan to co Diem (-10000)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row at position 18 (shifts rows 18.. down by one) -----
$ws.Rows("18:18").Insert()

# Copy the formatting (border/number-format) of the row below (which used to
# be row 18, now row 19) onto the freshly inserted blank row 18, then fill
# in the new H/I values.
$ws.Range("H19:I19").Copy()
$ws.Range("H18:I18").PasteSpecial(-4122)
$ws.Range("H18").Value = 45293
$ws.Range("I18").Value = 10

# --- 2) Update the (now shifted) last transfer-to-co-Diem row (was row 63,
#        is now row 64): 700k -> 1.7tr, amount 700 -> 1700 -------------------
$origText = $ws.Range("B64").Value()
$prefix = $origText.Substring(0, $origText.Length - 4)
$ws.Range("B64").Value = $prefix + "1.7tr"
$ws.Range("C64").Value = 1700

# --- 3) Append a brand-new row 65: Duy lends co Diem another 10tr ----------
$loanText = $ws.Range("B10").Value()
$ws.Range("B65").Value = $loanText
$ws.Range("C65").Value = -10000
$ws.Range("D65").Formula = "=D64+C65"

# --- 4) Update the cell cursor / selection to match where the user ended up
$ws.Range("D68").Select()
